$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1. Remove "as " before "this is a self-reported open forum."
Replace-Text "It is important to note that as this is a self-reported open forum." "It is important to note that this is a self-reported open forum."

# 2. Fix "inteval" -> "interval"
Replace-Text "building a confidence inteval for the average?" "building a confidence interval for the average?"

# 3. Fix "unimodel" -> "unimodal"
Replace-Text "Is it actually unimodel and symmetric?" "Is it actually unimodal and symmetric?"

# 4. Fix "there are any outliers" -> "there aren't any outliers"
Replace-Text "Luckily, there are any outliers" "Luckily, there aren’t any outliers"

# 5. Fix "such an Normal" -> "such as Normal"
Replace-Text "assessing tools such an Normal QQ plots" "assessing tools such as Normal QQ plots"

# 6. Insert comma: "above the average she" -> "above the average, she"
Replace-Text "Since Melissa’s ratio is above the average she would need" "Since Melissa’s ratio is above the average, she would need"

# 7. Remove "to " in "even more than to the amount"
Replace-Text "so that it is even more than to the amount she can flat dumbbell press" "so that it is even more than the amount she can flat dumbbell press"

# Remove the _GoBack bookmark (Word discards this automatically when the
# document is closed after an editing session)
try {
    $bm = $d.Bookmarks("_GoBack")
    if ($bm -ne $null) {
        $bm.Delete()
    }
} catch {
}
